# Auto-generated edit script applying the Sargatanas_Profits market-data refresh.
# Columns: H=currentAveragePrice I=currentAveragePriceNQ J=currentAveragePriceHQ
#          K=LevePriceNQ L=LevePriceHQ M=LeveProfitNQ N=LeveProfitHQ
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(8, 8).Value = 464.42554  # H8: 450.9375 -> 464.42554
$ws.Cells.Item(8, 9).Value = 1863.8  # I8: 1564 -> 1863.8
$ws.Cells.Item(8, 10).Value = 297.83334  # J8: 291.92856 -> 297.83334
$ws.Cells.Item(8, 11).Value = 5591.4  # K8: 4692 -> 5591.4
$ws.Cells.Item(8, 12).Value = 893.5000200000001  # L8: 875.78568 -> 893.5000200000001
$ws.Cells.Item(8, 13).Value = -5452.4  # M8: -4553 -> -5452.4
$ws.Cells.Item(8, 14).Value = -1171.50002  # N8: -1153.78568 -> -1171.50002
$ws.Cells.Item(53, 8).Value = 209.66667  # H53: 95 -> 209.66667
$ws.Cells.Item(53, 9).Value = 209.66667  # I53: 95 -> 209.66667
$ws.Cells.Item(53, 11).Value = 209.66667  # K53: 95 -> 209.66667
$ws.Cells.Item(53, 13).Value = 427.33333  # M53: 542 -> 427.33333
$ws.Cells.Item(86, 8).Value = 90813700  # H86: 90813710 -> 90813700
$ws.Cells.Item(86, 10).Value = 7938661  # J86: 7938696.5 -> 7938661
$ws.Cells.Item(86, 12).Value = 7938661  # L86: 7938696.5 -> 7938661
$ws.Cells.Item(86, 14).Value = -7940907  # N86: -7940942.5 -> -7940907
$ws.Cells.Item(89, 8).Value = 90813700  # H89: 90813710 -> 90813700
$ws.Cells.Item(89, 10).Value = 7938661  # J89: 7938696.5 -> 7938661
$ws.Cells.Item(89, 12).Value = 39693305  # L89: 39693482.5 -> 39693305
$ws.Cells.Item(89, 14).Value = -39704537  # N89: -39704714.5 -> -39704537
$ws.Cells.Item(123, 8).Value = 54999  # H123: 0 -> 54999
$ws.Cells.Item(123, 10).Value = 54999  # J123: 0 -> 54999
$ws.Cells.Item(123, 12).Value = 54999  # L123: 0 -> 54999
$ws.Cells.Item(123, 14).Value = -64799  # N123: None -> -64799
$ws.Cells.Item(125, 8).Value = 76924640  # H125: 90910696 -> 76924640
$ws.Cells.Item(125, 9).Value = 111112320  # I125: 142858320 -> 111112320
$ws.Cells.Item(125, 11).Value = 1000010880  # K125: 1285724880 -> 1000010880
$ws.Cells.Item(125, 13).Value = -1000008420  # M125: -1285722420 -> -1000008420
$ws.Cells.Item(132, 8).Value = 1671.0714  # H132: 2035.091 -> 1671.0714
$ws.Cells.Item(132, 9).Value = 1671.0714  # I132: 2035.091 -> 1671.0714
$ws.Cells.Item(132, 11).Value = 5013.2142  # K132: 6105.272999999999 -> 5013.2142
$ws.Cells.Item(132, 13).Value = -2483.2142  # M132: -3575.272999999999 -> -2483.2142
$ws.Cells.Item(137, 8).Value = 4575.243  # H137: 4303.45 -> 4575.243
$ws.Cells.Item(137, 9).Value = 2370.4583  # I137: 2212.7778 -> 2370.4583
$ws.Cells.Item(137, 11).Value = 7111.374899999999  # K137: 6638.3334 -> 7111.374899999999
$ws.Cells.Item(137, 13).Value = -4561.374899999999  # M137: -4088.3334 -> -4561.374899999999
$ws.Cells.Item(138, 8).Value = 1542182.8  # H138: 1670441 -> 1542182.8
$ws.Cells.Item(138, 9).Value = 1718.2142  # I138: 1667.9259 -> 1718.2142
$ws.Cells.Item(138, 10).Value = 2707939.5  # J138: 3035800.8 -> 2707939.5
$ws.Cells.Item(138, 11).Value = 5154.642599999999  # K138: 5003.7777 -> 5154.642599999999
$ws.Cells.Item(138, 12).Value = 8123818.5  # L138: 9107402.399999999 -> 8123818.5
$ws.Cells.Item(138, 13).Value = -14.64259999999922  # M138: 136.2223000000004 -> -14.64259999999922
$ws.Cells.Item(138, 14).Value = -8134098.5  # N138: -9117682.399999999 -> -8134098.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 4885482  # H32: 3453603.5 -> 4885482
$ws.Cells.Item(32, 9).Value = 5268402  # I32: 3640063.8 -> 5268402
$ws.Cells.Item(32, 11).Value = 5268402  # K32: 3640063.8 -> 5268402
$ws.Cells.Item(32, 13).Value = -5268115  # M32: -3639776.8 -> -5268115
$ws.Cells.Item(61, 8).Value = 8952.111000000001  # H61: 8957.666999999999 -> 8952.111000000001
$ws.Cells.Item(61, 9).Value = 2277.375  # I61: 2289.875 -> 2277.375
$ws.Cells.Item(61, 11).Value = 2277.375  # K61: 2289.875 -> 2277.375
$ws.Cells.Item(61, 13).Value = -2065.375  # M61: -2077.875 -> -2065.375
$ws.Cells.Item(121, 8).Value = 0  # H121: 59376 -> 0
$ws.Cells.Item(121, 10).Value = 0  # J121: 59376 -> 0
$ws.Cells.Item(121, 12).Value = 0  # L121: 59376 -> 0
$ws.Cells.Item(121, 14).ClearContents()  # N121: -62870 -> (removed)
$ws.Cells.Item(136, 8).Value = 8952.111000000001  # H136: 8957.666999999999 -> 8952.111000000001
$ws.Cells.Item(136, 9).Value = 2277.375  # I136: 2289.875 -> 2277.375
$ws.Cells.Item(136, 11).Value = 6832.125  # K136: 6869.625 -> 6832.125
$ws.Cells.Item(136, 13).Value = -4282.125  # M136: -4319.625 -> -4282.125

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(107, 8).Value = 45005096  # H107: 46880270 -> 45005096
$ws.Cells.Item(107, 9).Value = 48918320  # I107: 51141840 -> 48918320
$ws.Cells.Item(107, 11).Value = 48918320  # K107: 51141840 -> 48918320
$ws.Cells.Item(107, 13).Value = -48916400  # M107: -51139920 -> -48916400
$ws.Cells.Item(134, 8).Value = 4550.9  # H134: 4477.314 -> 4550.9
$ws.Cells.Item(134, 9).Value = 1652.8572  # I134: 1661.3715 -> 1652.8572
$ws.Cells.Item(134, 10).Value = 11313  # J134: 10637.1875 -> 11313
$ws.Cells.Item(134, 11).Value = 4958.571599999999  # K134: 4984.1145 -> 4958.571599999999
$ws.Cells.Item(134, 12).Value = 33939  # L134: 31911.5625 -> 33939
$ws.Cells.Item(134, 13).Value = -2423.571599999999  # M134: -2449.1145 -> -2423.571599999999
$ws.Cells.Item(134, 14).Value = -39009  # N134: -36981.5625 -> -39009

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 5743.8984  # H31: 6045.1113 -> 5743.8984
$ws.Cells.Item(31, 9).Value = 1525.0667  # I31: 1601 -> 1525.0667
$ws.Cells.Item(31, 10).Value = 10108.207  # J31: 10489.223 -> 10108.207
$ws.Cells.Item(31, 11).Value = 1525.0667  # K31: 1601 -> 1525.0667
$ws.Cells.Item(31, 12).Value = 10108.207  # L31: 10489.223 -> 10108.207
$ws.Cells.Item(31, 13).Value = -1230.0667  # M31: -1306 -> -1230.0667
$ws.Cells.Item(31, 14).Value = -10698.207  # N31: -11079.223 -> -10698.207
$ws.Cells.Item(34, 8).Value = 5743.8984  # H34: 6045.1113 -> 5743.8984
$ws.Cells.Item(34, 9).Value = 1525.0667  # I34: 1601 -> 1525.0667
$ws.Cells.Item(34, 10).Value = 10108.207  # J34: 10489.223 -> 10108.207
$ws.Cells.Item(34, 11).Value = 1525.0667  # K34: 1601 -> 1525.0667
$ws.Cells.Item(34, 12).Value = 10108.207  # L34: 10489.223 -> 10108.207
$ws.Cells.Item(34, 13).Value = -1323.0667  # M34: -1399 -> -1323.0667
$ws.Cells.Item(34, 14).Value = -10512.207  # N34: -10893.223 -> -10512.207
$ws.Cells.Item(58, 8).Value = 5189.3076  # H58: 5259.706 -> 5189.3076
$ws.Cells.Item(58, 10).Value = 7744.276  # J58: 7963.75 -> 7744.276
$ws.Cells.Item(58, 12).Value = 7744.276  # L58: 7963.75 -> 7744.276
$ws.Cells.Item(58, 14).Value = -8150.276  # N58: -8369.75 -> -8150.276
$ws.Cells.Item(93, 8).Value = 11872.833  # H93: 10177.286 -> 11872.833
$ws.Cells.Item(93, 9).Value = 3459  # I93: 2883.1667 -> 3459
$ws.Cells.Item(93, 11).Value = 3459  # K93: 2883.1667 -> 3459
$ws.Cells.Item(93, 13).Value = -1587  # M93: -1011.1667 -> -1587
$ws.Cells.Item(105, 8).Value = 4763020.5  # H105: 5103251.5 -> 4763020.5
$ws.Cells.Item(105, 9).Value = 6494029  # I105: 7143452.5 -> 6494029
$ws.Cells.Item(105, 11).Value = 6494029  # K105: 7143452.5 -> 6494029
$ws.Cells.Item(105, 13).Value = -6492282  # M105: -7141705.5 -> -6492282
$ws.Cells.Item(132, 8).Value = 6639.086  # H132: 6485.1943 -> 6639.086
$ws.Cells.Item(132, 9).Value = 4356.4736  # I132: 4193.6 -> 4356.4736
$ws.Cells.Item(132, 11).Value = 13069.4208  # K132: 12580.8 -> 13069.4208
$ws.Cells.Item(132, 13).Value = -10539.4208  # M132: -10050.8 -> -10539.4208
$ws.Cells.Item(134, 8).Value = 6581.231  # H134: 6124.852 -> 6581.231
$ws.Cells.Item(134, 9).Value = 1866.4166  # I134: 1489.7142 -> 1866.4166
$ws.Cells.Item(134, 10).Value = 10622.5  # J134: 11116.538 -> 10622.5
$ws.Cells.Item(134, 11).Value = 5599.2498  # K134: 4469.142599999999 -> 5599.2498
$ws.Cells.Item(134, 12).Value = 31867.5  # L134: 33349.614 -> 31867.5
$ws.Cells.Item(134, 13).Value = -3064.2498  # M134: -1934.142599999999 -> -3064.2498
$ws.Cells.Item(134, 14).Value = -36937.5  # N134: -38419.614 -> -36937.5
$ws.Cells.Item(136, 8).Value = 5189.3076  # H136: 5259.706 -> 5189.3076
$ws.Cells.Item(136, 10).Value = 7744.276  # J136: 7963.75 -> 7744.276
$ws.Cells.Item(136, 12).Value = 23232.828  # L136: 23891.25 -> 23232.828
$ws.Cells.Item(136, 14).Value = -28332.828  # N136: -28991.25 -> -28332.828
$ws.Cells.Item(141, 8).Value = 520178  # H141: 552999.7 -> 520178
$ws.Cells.Item(141, 10).Value = 520178  # J141: 552999.7 -> 520178
$ws.Cells.Item(141, 12).Value = 520178  # L141: 552999.7 -> 520178
$ws.Cells.Item(141, 14).Value = -530538  # N141: -563359.7 -> -530538

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(4, 8).Value = 3326755.5  # H4: 3659418.5 -> 3326755.5
$ws.Cells.Item(4, 9).Value = 4766877  # I4: 5561336 -> 4766877
$ws.Cells.Item(4, 11).Value = 14300631  # K4: 16684008 -> 14300631
$ws.Cells.Item(4, 13).Value = -14300519  # M4: -16683896 -> -14300519
$ws.Cells.Item(39, 8).Value = 8762.817999999999  # H39: 9035.091 -> 8762.817999999999
$ws.Cells.Item(39, 10).Value = 9913.429  # J39: 10341.286 -> 9913.429
$ws.Cells.Item(39, 12).Value = 29740.287  # L39: 31023.858 -> 29740.287
$ws.Cells.Item(39, 14).Value = -30328.287  # N39: -31611.858 -> -30328.287
$ws.Cells.Item(117, 8).Value = 1200  # H117: 1450 -> 1200
$ws.Cells.Item(117, 10).Value = 0  # J117: 1950 -> 0
$ws.Cells.Item(117, 12).Value = 0  # L117: 5850 -> 0
$ws.Cells.Item(117, 14).ClearContents()  # N117: -12734 -> (removed)

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(136, 8).Value = 25214.438  # H136: 27515.068 -> 25214.438
$ws.Cells.Item(136, 9).Value = 25800  # I136: 30000 -> 25800
$ws.Cells.Item(136, 10).Value = 25106  # J136: 27117.48 -> 25106
$ws.Cells.Item(136, 11).Value = 77400  # K136: 90000 -> 77400
$ws.Cells.Item(136, 12).Value = 75318  # L136: 81352.44 -> 75318
$ws.Cells.Item(136, 13).Value = -74850  # M136: -87450 -> -74850
$ws.Cells.Item(136, 14).Value = -80418  # N136: -86452.44 -> -80418

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(22, 8).Value = 1091.6471  # H22: 1234 -> 1091.6471
$ws.Cells.Item(22, 9).Value = 404.85715  # I22: 439.58334 -> 404.85715
$ws.Cells.Item(22, 10).Value = 4296.6665  # J22: 6000.5 -> 4296.6665
$ws.Cells.Item(22, 11).Value = 404.85715  # K22: 439.58334 -> 404.85715
$ws.Cells.Item(22, 12).Value = 4296.6665  # L22: 6000.5 -> 4296.6665
$ws.Cells.Item(22, 13).Value = -109.85715  # M22: -144.58334 -> -109.85715
$ws.Cells.Item(22, 14).Value = -4886.6665  # N22: -6590.5 -> -4886.6665
$ws.Cells.Item(27, 8).Value = 1091.6471  # H27: 1234 -> 1091.6471
$ws.Cells.Item(27, 9).Value = 404.85715  # I27: 439.58334 -> 404.85715
$ws.Cells.Item(27, 10).Value = 4296.6665  # J27: 6000.5 -> 4296.6665
$ws.Cells.Item(27, 11).Value = 404.85715  # K27: 439.58334 -> 404.85715
$ws.Cells.Item(27, 12).Value = 4296.6665  # L27: 6000.5 -> 4296.6665
$ws.Cells.Item(27, 13).Value = -297.85715  # M27: -332.58334 -> -297.85715
$ws.Cells.Item(27, 14).Value = -4510.6665  # N27: -6214.5 -> -4510.6665
$ws.Cells.Item(40, 8).Value = 5611.077  # H40: 5523.1665 -> 5611.077
$ws.Cells.Item(40, 10).Value = 7777.5  # J40: 7999.8 -> 7777.5
$ws.Cells.Item(40, 12).Value = 7777.5  # L40: 7999.8 -> 7777.5
$ws.Cells.Item(40, 14).Value = -8049.5  # N40: -8271.799999999999 -> -8049.5
$ws.Cells.Item(61, 8).Value = 4592.7607  # H61: 4739.636 -> 4592.7607
$ws.Cells.Item(61, 9).Value = 3278.0938  # I61: 3360.516 -> 3278.0938
$ws.Cells.Item(61, 10).Value = 7597.7144  # J61: 8028.3076 -> 7597.7144
$ws.Cells.Item(61, 11).Value = 3278.0938  # K61: 3360.516 -> 3278.0938
$ws.Cells.Item(61, 12).Value = 7597.7144  # L61: 8028.3076 -> 7597.7144
$ws.Cells.Item(61, 13).Value = -3076.0938  # M61: -3158.516 -> -3076.0938
$ws.Cells.Item(61, 14).Value = -8001.7144  # N61: -8432.3076 -> -8001.7144
$ws.Cells.Item(93, 8).Value = 886  # H93: 686.5 -> 886
$ws.Cells.Item(93, 9).Value = 545  # I93: 545.625 -> 545
$ws.Cells.Item(93, 10).Value = 2250  # J93: 1250 -> 2250
$ws.Cells.Item(93, 11).Value = 545  # K93: 545.625 -> 545
$ws.Cells.Item(93, 12).Value = 2250  # L93: 1250 -> 2250
$ws.Cells.Item(93, 13).Value = 703  # M93: 702.375 -> 703
$ws.Cells.Item(93, 14).Value = -4746  # N93: -3746 -> -4746
$ws.Cells.Item(113, 8).Value = 4592.7607  # H113: 4739.636 -> 4592.7607
$ws.Cells.Item(113, 9).Value = 3278.0938  # I113: 3360.516 -> 3278.0938
$ws.Cells.Item(113, 10).Value = 7597.7144  # J113: 8028.3076 -> 7597.7144
$ws.Cells.Item(113, 11).Value = 3278.0938  # K113: 3360.516 -> 3278.0938
$ws.Cells.Item(113, 12).Value = 7597.7144  # L113: 8028.3076 -> 7597.7144
$ws.Cells.Item(113, 13).Value = -1108.0938  # M113: -1190.516 -> -1108.0938
$ws.Cells.Item(113, 14).Value = -11937.7144  # N113: -12368.3076 -> -11937.7144
$ws.Cells.Item(122, 8).Value = 3919  # H122: 4086 -> 3919
$ws.Cells.Item(122, 9).Value = 2508  # I122: 2609.9 -> 2508
$ws.Cells.Item(122, 11).Value = 7524  # K122: 7829.700000000001 -> 7524
$ws.Cells.Item(122, 13).Value = -5074  # M122: -5379.700000000001 -> -5074
$ws.Cells.Item(132, 8).Value = 6278.9214  # H132: 6348.5 -> 6278.9214
$ws.Cells.Item(132, 9).Value = 3314.1304  # I132: 3337.5 -> 3314.1304
$ws.Cells.Item(132, 11).Value = 9942.3912  # K132: 10012.5 -> 9942.3912
$ws.Cells.Item(132, 13).Value = -7412.3912  # M132: -7482.5 -> -7412.3912
$ws.Cells.Item(136, 8).Value = 10617.456  # H136: 11092.463 -> 10617.456
$ws.Cells.Item(136, 9).Value = 2824.0344  # I136: 2911.3462 -> 2824.0344
$ws.Cells.Item(136, 11).Value = 8472.1032  # K136: 8734.0386 -> 8472.1032
$ws.Cells.Item(136, 13).Value = -5922.1032  # M136: -6184.0386 -> -5922.1032

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(101, 8).Value = 24971  # H101: 24256.5 -> 24971
$ws.Cells.Item(101, 10).Value = 24971  # J101: 24256.5 -> 24971
$ws.Cells.Item(101, 12).Value = 24971  # L101: 24256.5 -> 24971
$ws.Cells.Item(101, 14).Value = -31461  # N101: -30746.5 -> -31461
$ws.Cells.Item(122, 8).Value = 2651.825  # H122: 2762.5898 -> 2651.825
$ws.Cells.Item(122, 9).Value = 2084.5715  # I122: 2172.5356 -> 2084.5715
$ws.Cells.Item(122, 10).Value = 3975.4167  # J122: 4264.5454 -> 3975.4167
$ws.Cells.Item(122, 11).Value = 6253.7145  # K122: 6517.6068 -> 6253.7145
$ws.Cells.Item(122, 12).Value = 11926.2501  # L122: 12793.6362 -> 11926.2501
$ws.Cells.Item(122, 13).Value = -3803.7145  # M122: -4067.6068 -> -3803.7145
$ws.Cells.Item(122, 14).Value = -16826.2501  # N122: -17693.6362 -> -16826.2501
$ws.Cells.Item(132, 8).Value = 4688.6577  # H132: 4489.225 -> 4688.6577
$ws.Cells.Item(132, 9).Value = 4338.8335  # I132: 4111.4062 -> 4338.8335
$ws.Cells.Item(132, 11).Value = 13016.5005  # K132: 12334.2186 -> 13016.5005
$ws.Cells.Item(132, 13).Value = -10486.5005  # M132: -9804.2186 -> -10486.5005
$ws.Cells.Item(136, 8).Value = 2547.9143  # H136: 3789.9285 -> 2547.9143
$ws.Cells.Item(136, 10).Value = 4352.385  # J136: 6329.05 -> 4352.385
$ws.Cells.Item(136, 12).Value = 13057.155  # L136: 18987.15 -> 13057.155
$ws.Cells.Item(136, 14).Value = -18157.155  # N136: -24087.15 -> -18157.155

